# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" worksheet right after "总计" (position 2), pushing
#    all the existing quarter sheets one slot later - done by copying the
#    existing "2022-Q2" sheet (same layout/styling) and renaming it.
# 2) Populate it with the fund-holdings detail rows for 2022-Q3.
# 3) Update the "总计" (summary) sheet: insert a new row right under the
#    header with the 2022-Q3 totals, pushing the previously-existing rows
#    down by one and renumbering the index column (A).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new sheet by copying the existing "2022-Q2" sheet (so it
# inherits the same column widths / header & index-column styling), then
# rename it and move it into the right slot (right after "总计").
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$templateSheet.Copy([System.Reflection.Missing]::Value, $totalSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template has 9 data rows (2..10); we need 21 (2..22). Extend the
# index-column (A) styling down to row 22 before writing values.
$newSheet.Range("A2").Copy()
$newSheet.Range("A11:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fund-code/name/size/position columns are stored as text in this workbook
# (not numbers) - format as Text first so values like "005233" keep their
# leading zero and "42.36" etc. aren't reinterpreted as numbers.
$newSheet.Range("B2:G22").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Step 2: fill in the header + 21 data rows for 2022-Q3.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "005233", "广发睿毅领先混合A",            "42.36", "89.88", "4.95", "2.0968", 8),
    @(1,  "001071", "华安媒体互联网混合A",           "48.36", "89.21", "3.44", "1.6636", 5),
    @(2,  "001694", "华安沪港深外延增长混合A",        "38.96", "89.76", "3.90", "1.5194", 6),
    @(3,  "001763", "广发多策略灵活配置混合",         "22.52", "90.50", "4.56", "1.0269", 9),
    @(4,  "006879", "华安智能生活混合A",             "27.26", "88.28", "3.74", "1.0195", 5),
    @(5,  "012449", "广发睿毅领先混合C",             "19.65", "89.88", "4.95", "0.9727", 8),
    @(6,  "012528", "广发鑫睿一年持有期混合A",        "7.26",  "92.69", "4.60", "0.3340", 10),
    @(7,  "014734", "广发睿合混合A",                "5.96",  "86.96", "4.65", "0.2771", 10),
    @(8,  "013621", "华安智能生活混合C",             "6.84",  "88.28", "3.74", "0.2558", 5),
    @(9,  "000880", "富国研究精选灵活配置混合A",       "3.84",  "89.44", "6.24", "0.2396", 5),
    @(10, "012529", "广发鑫睿一年持有期混合C",        "4.74",  "92.69", "4.60", "0.2180", 10),
    @(11, "014754", "华安景气优选混合A",             "5.37",  "87.64", "3.45", "0.1853", 5),
    @(12, "014177", "华安景气驱动一年持有混合A",       "3.38",  "89.22", "3.30", "0.1115", 6),
    @(13, "014735", "广发睿合混合C",                "1.47",  "86.96", "4.65", "0.0684", 10),
    @(14, "014755", "华安景气优选混合C",             "1.53",  "87.64", "3.45", "0.0528", 5),
    @(15, "013620", "华安媒体互联网混合C",           "1.23",  "89.21", "3.44", "0.0423", 5),
    @(16, "014178", "华安景气驱动一年持有混合C",       "0.32",  "89.22", "3.30", "0.0106", 6),
    @(17, "016313", "富国研究精选灵活配置混合C",       "0.10",  "89.44", "6.24", "0.0062", 5),
    @(18, "014972", "华安沪港深外延增长混合C",        "0.07",  "89.76", "3.90", "0.0027", 6),
    @(19, "001914", "中信建投聚利混合A",             "0.10",  "39.73", "2.04", "0.0020", 8),
    @(20, "006845", "中信建投聚利混合C",             "0.01",  "39.73", "2.04", "0.0002", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" sheet - shift existing rows 2-8 down to 3-9 and
# write the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------------

# Row 9 is brand new - give its index cell (A9) the same style as A8 (bold,
# bordered, centered) before writing the value into it.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($row = 8; $row -ge 2; $row--) {
    $newRow = $row + 1
    $totalSheet.Cells.Item($newRow, 1).Value = $row - 1
    $totalSheet.Cells.Item($newRow, 2).Value = $totalSheet.Cells.Item($row, 2).Value2
    $totalSheet.Cells.Item($newRow, 3).Value = $totalSheet.Cells.Item($row, 3).Value2
    $totalSheet.Cells.Item($newRow, 4).Value = $totalSheet.Cells.Item($row, 4).Value2
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 21
$totalSheet.Cells.Item(2, 4).Value = 10.11
